$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = 'D2'; Value = '34.081.21' }
    @{ Cell = 'E2'; Value = '  -1.51%  ' }
    @{ Cell = 'D3'; Value = '1.794.75' }
    @{ Cell = 'E3'; Value = '  -1.84%  ' }
    @{ Cell = 'E4'; Value = '  +0.46%  ' }
    @{ Cell = 'D5'; Value = '228.29' }
    @{ Cell = 'E5'; Value = '  -3.13%  ' }
    @{ Cell = 'D6'; Value = '0.556' }
    @{ Cell = 'E6'; Value = '  +0.49%  ' }
    @{ Cell = 'E7'; Value = '  +0.52%  ' }
    @{ Cell = 'D8'; Value = '31.27' }
    @{ Cell = 'E8'; Value = '  -3.00%  ' }
    @{ Cell = 'D9'; Value = '46.12' }
    @{ Cell = 'E9'; Value = '  -0.48%  ' }
    @{ Cell = 'E10'; Value = '  -1.95%  ' }
    @{ Cell = 'D11'; Value = '0.0663' }
    @{ Cell = 'E11'; Value = '  -3.52%  ' }
    @{ Cell = 'E12'; Value = '  -0.19%  ' }
    @{ Cell = 'D13'; Value = '2.050.59' }
    @{ Cell = 'E13'; Value = '  -1.73%  ' }
    @{ Cell = 'D14'; Value = '11.31' }
    @{ Cell = 'E14'; Value = '  +8.80%  ' }
    @{ Cell = 'D15'; Value = '1.793.23' }
    @{ Cell = 'E15'; Value = '  -1.84%  ' }
    @{ Cell = 'D16'; Value = '0.636' }
    @{ Cell = 'E16'; Value = '  -2.38%  ' }
    @{ Cell = 'D17'; Value = '34.084.56' }
    @{ Cell = 'E17'; Value = '  -1.40%  ' }
    @{ Cell = 'D18'; Value = '4.23' }
    @{ Cell = 'E18'; Value = '  -3.72%  ' }
    @{ Cell = 'D19'; Value = '69.79' }
    @{ Cell = 'E19'; Value = '  -3.08%  ' }
    @{ Cell = 'D20'; Value = '253.89' }
    @{ Cell = 'E20'; Value = '  -5.00%  ' }
    @{ Cell = 'E21'; Value = '  -2.56%  ' }
    @{ Cell = 'E22'; Value = '  +0.27%  ' }
    @{ Cell = 'D23'; Value = '10.47' }
    @{ Cell = 'E23'; Value = '  -1.63%  ' }
    @{ Cell = 'D24'; Value = '4.30' }
    @{ Cell = 'E24'; Value = '  -3.74%  ' }
    @{ Cell = 'E25'; Value = '  -1.90%  ' }
    @{ Cell = 'D26'; Value = '157.73' }
    @{ Cell = 'E26'; Value = '  -2.94%  ' }
    @{ Cell = 'D27'; Value = '16.65' }
    @{ Cell = 'E27'; Value = '  -3.28%  ' }
    @{ Cell = 'D28'; Value = '7.04' }
    @{ Cell = 'E28'; Value = '  -2.42%  ' }
    @{ Cell = 'D29'; Value = '0.115' }
    @{ Cell = 'E29'; Value = '  -2.57%  ' }
    @{ Cell = 'E30'; Value = '  +0.68%  ' }
    @{ Cell = 'E31'; Value = '  +0.24%  ' }
    @{ Cell = 'E32'; Value = '  -0.42%  ' }
    @{ Cell = 'E33'; Value = '  -1.06%  ' }
    @{ Cell = 'E34'; Value = '  +0.63%  ' }
    @{ Cell = 'E35'; Value = '  -0.57%  ' }
    @{ Cell = 'D36'; Value = '1.488.30' }
    @{ Cell = 'E36'; Value = '  -7.06%  ' }
    @{ Cell = 'E37'; Value = '  -0.46%  ' }
    @{ Cell = 'E38'; Value = '  +0.23%  ' }
    @{ Cell = 'D39'; Value = '0.0188' }
    @{ Cell = 'E39'; Value = '  -1.18%  ' }
    @{ Cell = 'D40'; Value = '84.06' }
    @{ Cell = 'E40'; Value = '  -6.57%  ' }
    @{ Cell = 'D41'; Value = '2.84' }
    @{ Cell = 'E41'; Value = '  -0.79%  ' }
    @{ Cell = 'E42'; Value = '  -0.38%  ' }
    @{ Cell = 'D43'; Value = '0.907' }
    @{ Cell = 'E43'; Value = '  -3.33%  ' }
    @{ Cell = 'E44'; Value = '  -4.61%  ' }
    @{ Cell = 'D45'; Value = '0.0515' }
    @{ Cell = 'E45'; Value = '  -1.21%  ' }
    @{ Cell = 'E46'; Value = '  +1.84%  ' }
    @{ Cell = 'D47'; Value = '1.949.26' }
    @{ Cell = 'E47'; Value = '  -1.09%  ' }
    @{ Cell = 'E48'; Value = '  -1.64%  ' }
    @{ Cell = 'E49'; Value = '  +0.20%  ' }
    @{ Cell = 'D50'; Value = '11.84' }
    @{ Cell = 'E50'; Value = '  +1.87%  ' }
    @{ Cell = 'D51'; Value = '51.62' }
    @{ Cell = 'E51'; Value = '  -5.46%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
